# Automatizada la carga de los archivos por jornada al archivo principal.
# Minutos, Goles, Asistencias
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 21

# Add new headers "Goles" (AP1) and "Asistencias de Gol" (AQ1), copying the
# existing header style/format from AO1 ("Capitan") so the new columns look
# consistent with the rest of the header row.
$ws.Range("AO1").Copy()
$ws.Range("AP1:AQ1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("AP1").Value = "Goles"
$ws.Range("AQ1").Value = "Asistencias de Gol"

# Clear out the old "Capitan" (Si/No) values in column AO and populate the
# new "Goles" / "Asistencias de Gol" columns with 0 for every player row.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 41).Value = ""
    $ws.Cells.Item($r, 42).Value = 0
    $ws.Cells.Item($r, 43).Value = 0
}

# Fix the accent on "Si" -> "Sí" for the substitute players (rows 13-21).
for ($r = 13; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 6).Value = "Sí"
}
